# Edit script: update vm_pu.xlsx values per commit "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns in order matching the data rows below: B, C, D, E, F, I, J, K, L, M
$cols = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M")

# Each inner array holds the new values for row 2..25, in the column order above.
$data = @(
    @(1.02, 1.075347494238983, 1.074379477558459, 1.078507792967299, 1.087576891009372, 1.046326922860092, 1.080252411879709, 1.077068845344237, 1.081186282913323, 1.090231797435664),
    @(1.02, 1.077035760507893, 1.075711047986396, 1.079997171536834, 1.089111649988354, 1.046708783593401, 1.081596655948686, 1.078216266905025, 1.082491921752737, 1.091584432683649),
    @(1.02, 1.07812636344143, 1.076570790845498, 1.080959382884945, 1.090103217238498, 1.046953668928187, 1.08246427778631, 1.078956307803063, 1.083334725729606, 1.092457643699512),
    @(1.02, 1.07858442836955, 1.076931786110253, 1.081363542114097, 1.090519715901062, 1.047056093813832, 1.08282850911782, 1.079266849327885, 1.08368856206923, 1.09282426265734),
    @(1.02, 1.078661314901923, 1.076992373184148, 1.081431381691099, 1.090589627223858, 1.047073260701551, 1.082889635092003, 1.079318957298889, 1.083747944920359, 1.092885791741686),
    @(1.02, 1.078132485780691, 1.076575616202676, 1.080964784661053, 1.090108783901477, 1.046955039593173, 1.082469146682579, 1.078960459510462, 1.083339455577029, 1.092462544350638),
    @(1.02, 1.07591843469287, 1.074829880262404, 1.079011454528691, 1.088095890916568, 1.046456432520007, 1.080707164773381, 1.077457126699344, 1.081627955129657, 1.090689352983394),
    @(1.02, 1.072002585176689, 1.071738989729476, 1.075557440223384, 1.084536847490893, 1.045560818173205, 1.077585148644582, 1.074789226581031, 1.078596155686682, 1.087548817532606),
    @(1.02, 1.069381653465449, 1.069668045167711, 1.073246157110562, 1.082155484694597, 1.044952133548824, 1.075491726565734, 1.072997512694851, 1.076563752887089, 1.085443894222537),
    @(1.02, 1.068244151090147, 1.06876874052649, 1.072243183490083, 1.081122148797083, 1.044685773592074, 1.074582268767457, 1.072218469546641, 1.075680930490794, 1.084529660186259),
    @(1.02, 1.067821225061398, 1.068434302959436, 1.071870297480896, 1.080737982049439, 1.044586412189103, 1.074243995829714, 1.071928605804064, 1.075352583735224, 1.084189644070814),
    @(1.02, 1.067911962759787, 1.068506059014197, 1.071950298213142, 1.080820402670567, 1.044607744765544, 1.074316577496773, 1.071990804986172, 1.075423034718709, 1.084262598250924),
    @(1.02, 1.068209200271689, 1.06874110394122, 1.072212367551063, 1.081090400447948, 1.044677569006026, 1.07455431646079, 1.072194519415773, 1.075653798017287, 1.084501563155275),
    @(1.02, 1.068392283786811, 1.068885870151444, 1.072373792154338, 1.081256709641429, 1.04472053380894, 1.074700734043716, 1.072319969083383, 1.075795921998691, 1.08464874018802),
    @(1.02, 1.069457089398969, 1.069727674068462, 1.073312674493368, 1.082224016736908, 1.044969751789146, 1.075552020367404, 1.073049146650222, 1.076622283529602, 1.085504509328316),
    @(1.02, 1.070124302707439, 1.070255020245219, 1.073901021472304, 1.082830189013911, 1.045125328775189, 1.076085201463149, 1.073505672425775, 1.077139887519965, 1.086040557169958),
    @(1.02, 1.070513224883639, 1.070562364827341, 1.074243985384741, 1.083183548466898, 1.045215804709953, 1.076395908968274, 1.07377164617727, 1.077441529742711, 1.086352956449703),
    @(1.02, 1.070645794936822, 1.07066711968548, 1.074360892094787, 1.083303999514134, 1.045246609059738, 1.076501803656188, 1.073862284034348, 1.077544336878928, 1.086459431312994),
    @(1.02, 1.070052743200398, 1.070198466661256, 1.073837919095829, 1.082765174324465, 1.045108664736536, 1.076028026041801, 1.07345672365985, 1.077084381242275, 1.085983072139223),
    @(1.02, 1.068121682574053, 1.068671900069764, 1.07213520404931, 1.081010902298604, 1.044657019224952, 1.074484321037545, 1.072134544258118, 1.0755858558494, 1.084431205866922),
    @(1.02, 1.066905186755777, 1.067709792487276, 1.071062684757618, 1.079905950555697, 1.044370600043724, 1.073511066912109, 1.071300384131476, 1.074641196853811, 1.083452999528733),
    @(1.02, 1.06755030244752, 1.068220044570811, 1.071631436350124, 1.080491896850685, 1.044522669830758, 1.074027263451748, 1.071742861624718, 1.075142216495939, 1.083971804252531),
    @(1.02, 1.07008507864028, 1.070224021556484, 1.073866432989238, 1.082794552314911, 1.045116195331064, 1.076053862043631, 1.073478842458202, 1.077109462966353, 1.086009047979511),
    @(1.02, 1.073016703076373, 1.072539847956031, 1.076451862814411, 1.085458432361665, 1.04579438906082, 1.078394354501358, 1.075481220890121, 1.079381886122286, 1.088362663312961),
)


for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowVals[$j]
    }
}
